# Update cryptocurrency price / volume(1h) data in the active worksheet,
# reflecting the latest scrape (commit: "Updated cryptos list ... with GitHub Actions").
# Rows 48/49 swap content (Mantle <-> VeChain) as their relative ranking changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.322.52'
$ws.Range("E2").Value = '  +2.90%  '
$ws.Range("D3").Value = '2.647.90'
$ws.Range("E3").Value = '  +1.84%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'603.81"
$ws.Range("E5").Value = '  +2.55%  '
$ws.Range("D6").Value = "'156.90"
$ws.Range("E6").Value = '  +5.13%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("E9").Value = '  +10.85%  '
$ws.Range("E10").Value = '  +6.05%  '
$ws.Range("D11").Value = "'5.82"
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("E12").Value = '  +2.72%  '
$ws.Range("D13").Value = "'29.42"
$ws.Range("E13").Value = '  +6.67%  '
$ws.Range("D14").Value = "'0.0000189"
$ws.Range("E14").Value = '  +21.57%  '
$ws.Range("D15").Value = '3.123.41'
$ws.Range("E15").Value = '  +1.95%  '
$ws.Range("D16").Value = '65.136.57'
$ws.Range("E16").Value = '  +2.91%  '
$ws.Range("D17").Value = '2.654.78'
$ws.Range("E17").Value = '  +2.74%  '
$ws.Range("E18").Value = '  +5.48%  '
$ws.Range("D19").Value = "'4.93"
$ws.Range("E19").Value = '  +4.23%  '
$ws.Range("D20").Value = "'360.19"
$ws.Range("E20").Value = '  +4.58%  '
$ws.Range("D21").Value = "'7.39"
$ws.Range("E21").Value = '  +8.22%  '
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = '  -0.12%  '
$ws.Range("D23").Value = "'69.42"
$ws.Range("D24").Value = "'1.72"
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("D25").Value = "'9.46"
$ws.Range("E25").Value = '  +3.07%  '
$ws.Range("E26").Value = '  +1.09%  '
$ws.Range("D27").Value = "'8.33"
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("E28").Value = '  +3.02%  '
$ws.Range("D29").Value = '0.0₃0975'
$ws.Range("E29").Value = '  +14.04%  '
$ws.Range("D30").Value = "'552.84"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").Value = "'2.21"
$ws.Range("E31").Value = '  +9.15%  '
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("D33").Value = "'1.81"
$ws.Range("E33").Value = '  +2.97%  '
$ws.Range("D34").Value = "'5.63"
$ws.Range("E34").Value = '  +6.68%  '
$ws.Range("E35").Value = '  +5.94%  '
$ws.Range("E36").Value = '  +4.90%  '
$ws.Range("D37").Value = "'20.57"
$ws.Range("E37").Value = '  +6.13%  '
$ws.Range("E38").Value = '  +3.77%  '
$ws.Range("D39").Value = "'162.25"
$ws.Range("E39").Value = '  -1.67%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = "'42.77"
$ws.Range("E42").Value = '  +8.26%  '
$ws.Range("D43").Value = "'166.90"
$ws.Range("E43").Value = '  +1.36%  '
$ws.Range("D44").Value = "'4.19"
$ws.Range("E44").Value = '  +3.81%  '
$ws.Range("D45").Value = "'0.0624"
$ws.Range("E45").Value = '  +7.84%  '
$ws.Range("E46").Value = '  +9.65%  '
$ws.Range("D47").Value = "'23.35"
$ws.Range("E47").Value = '  +2.18%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = "'0.0265"
$ws.Range("E48").Value = '  +6.42%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = "'0.655"
$ws.Range("E49").Value = '  +3.78%  '
$ws.Range("D50").Value = "'0.0983"
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("D51").Value = "'19.76"
$ws.Range("E51").Value = '  +3.81%  '
